# Actualización automática 2025-06-30 14:40:09
$wb = $excel.ActiveWorkbook

$wsVentasPorGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsVentaMensual    = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento    = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# --- Sheet "VENTAS POR GRUPO" ---
# Fila 6: CERAMICAS AL COSTO S.A.S. - PORCELANATO
$wsVentasPorGrupo.Range("M6").Value = 25364.28

# Fila 18: ZAMBRANO ANGELA MARIA - FREGADEROS DE COCINA / GRIFERIAS
$wsVentasPorGrupo.Range("E18").Value = 194.46
$wsVentasPorGrupo.Range("G18").Value = 122.22

# Fila 19: contadores "X de 17"
$wsVentasPorGrupo.Range("E19").Value = "2 de 17"
$wsVentasPorGrupo.Range("G19").Value = "2 de 17"

# --- Sheet "VENTA MENSUAL" ---
# Fila 6: CERAMICAS AL COSTO S.A.S. - junio
$wsVentaMensual.Range("F6").Value = 25364.28

# Fila 18: ZAMBRANO ANGELA MARIA - junio
$wsVentaMensual.Range("F18").Value = 4798.25

# Fila 19: TOTAL - junio
$wsVentaMensual.Range("F19").Value = 34453.24

# --- Sheet "CUMPLIMIENTO MENSUAL" ---
# Ensanchar columna E (POR CUMPLIR) de 23 a 24
$wsCumplimiento.Columns.Item(5).ColumnWidth = 23.15

# Fila 4: FREGADEROS DE COCINA
$wsCumplimiento.Range("D4").Value = 259.28
$wsCumplimiento.Range("E4").Value = 113.713863046034
$wsCumplimiento.Range("F4").Value = 0.6951320804117366

# Fila 6: GRIFERIAS
$wsCumplimiento.Range("D6").Value = 162.96
$wsCumplimiento.Range("E6").Value = -56.14000000000001
$wsCumplimiento.Range("F6").Value = 1.525557011795544

# Fila 16: PORCELANATO
$wsCumplimiento.Range("D16").Value = 31134.24
$wsCumplimiento.Range("E16").Value = -2924.400000000001
$wsCumplimiento.Range("F16").Value = 1.103665954858305

# Fila 19: TOTAL
$wsCumplimiento.Range("D19").Value = 34453.24000000001
$wsCumplimiento.Range("E19").Value = 12766.06386304603
$wsCumplimiento.Range("F19").Value = 0.7296431158732777
